$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data (columns A:K, rows 1:20) one column to the right
# (to B:L), so that a new column A can be inserted for the segment index.
# Work from the rightmost column to the leftmost so we never overwrite data
# before it has been copied. Range.Copy(Destination) copies both the value
# and the cell formatting/style in one shot.
$srcCols = @("K","J","I","H","G","F","E","D","C","B","A")
$dstCols = @("L","K","J","I","H","G","F","E","D","C","B")

for ($i = 0; $i -lt $srcCols.Length; $i++) {
    $src = $ws.Range($srcCols[$i] + "1:" + $srcCols[$i] + "20")
    $dst = $ws.Range($dstCols[$i] + "1:" + $dstCols[$i] + "20")
    $src.Copy($dst)
}

# New column A header.
$ws.Range("A1").Value = "segments"

# New column A values: a 0-based segment index for each data row.
for ($i = 0; $i -le 18; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
}

# Apply the same formatting used by the header row / index column to the
# new column A cells (bold, thin border, centered horizontally, top
# vertically) by copying the format from the now-shifted header/segment
# cells in column B.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("B2:B20").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)

$excel.CutCopyMode = 0
